$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newestDate = 44261.52689609415
$midDate    = 44261.50553879629
$oldestDate = 44261.48421208333

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newestDate
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $midDate
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $oldestDate
}
